$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Kommentar des Mentors" column (G) with header and per-row comments
$ws.Range("G1").Value = "Kommentar des Mentors"
$ws.Range("G1").Font.Bold = $true
$ws.Range("G2").Value = "Hallo Burak, deine Abgabe war sehr gut!"
$ws.Range("G3").Value = "Hallo Lionel, leider war deine Abgabe eine Katastrophe"
$ws.Range("G4").Value = "Hallo Cristiano, deine Abgabe war ganz okay"
$ws.Range("G5").Value = "Hallo lieber Anderson Talisca, sehr schönes Freistoßtor!"

# Widen the new column so the comments are fully visible
$ws.Columns.Item(7).ColumnWidth = 30.498697916666668

# Move/record the active selection as it was left after the edit
$ws.Range("G22").Select() | Out-Null
